# Scheduled market-data refresh: update the fetched average-price columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ) and
# the derived Leve price/profit columns for each class's profit sheet.
# All touched cells hold literal numbers (no formulas in this workbook), so
# each changed value is written directly; a few rows gain/lose a trailing
# profit cell depending on whether the computed value exists this run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 204.2
$ws.Range("I2").Value = 167.75
$ws.Range("K2").Value = 167.75
$ws.Range("M2").Value = -54.75

$ws.Range("H8").Value = 1177.5555
$ws.Range("I8").Value = 66.333336
$ws.Range("J8").Value = 3400
$ws.Range("K8").Value = 199.000008
$ws.Range("L8").Value = 10200
$ws.Range("M8").Value = -60.00000800000001
$ws.Range("N8").Value = -10478

$ws.Range("H15").Value = 184.75
$ws.Range("I15").Value = 184.75
$ws.Range("K15").Value = 554.25
$ws.Range("M15").Value = -385.25

$ws.Range("H29").Value = 521.4286
$ws.Range("J29").Value = 437.5
$ws.Range("L29").Value = 1312.5
$ws.Range("N29").Value = -1874.5

$ws.Range("H38").Value = 205
$ws.Range("I38").Value = 50.909092
$ws.Range("K38").Value = 152.727276
$ws.Range("M38").Value = 219.272724

$ws.Range("H43").Value = 5966
$ws.Range("I43").Value = 4966.6665
$ws.Range("J43").Value = 6965.3335
$ws.Range("K43").Value = 4966.6665
$ws.Range("L43").Value = 6965.3335
$ws.Range("M43").Value = -4897.6665
$ws.Range("N43").Value = -7103.3335

$ws.Range("H53").Value = 760.58826
$ws.Range("J53").Value = 806.2857
$ws.Range("L53").Value = 806.2857
$ws.Range("N53").Value = -2080.2857

$ws.Range("H74").Value = 7000
$ws.Range("I74").Value = 7000
$ws.Range("K74").Value = 7000
$ws.Range("M74").Value = -6064

$ws.Range("H77").Value = 7000
$ws.Range("I77").Value = 7000
$ws.Range("K77").Value = 35000
$ws.Range("M77").Value = -30320

$ws.Range("H98").Value = 1040.6666
$ws.Range("I98").Value = 1040.6666
$ws.Range("K98").Value = 1040.6666
$ws.Range("M98").Value = 457.3334

$ws.Range("H111").Value = 1990
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").Value = $null

$ws.Range("H122").Value = 1040.6666
$ws.Range("I122").Value = 1040.6666
$ws.Range("K122").Value = 3121.9998
$ws.Range("M122").Value = -671.9998000000001

$ws.Range("H132").Value = 3274.2
$ws.Range("I132").Value = 3092.75
$ws.Range("K132").Value = 9278.25
$ws.Range("M132").Value = -6748.25

$ws.Range("H137").Value = 827.9375
$ws.Range("I137").Value = 684.9
$ws.Range("J137").Value = 1066.3334
$ws.Range("K137").Value = 2054.7
$ws.Range("L137").Value = 3199.0002
$ws.Range("M137").Value = 495.3000000000002
$ws.Range("N137").Value = -8299.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3096.8572
$ws.Range("I2").Value = 1056.4445
$ws.Range("K2").Value = 1056.4445
$ws.Range("M2").Value = -943.4445000000001

$ws.Range("H74").Value = 907.1667
$ws.Range("I74").Value = 907.1667
$ws.Range("K74").Value = 907.1667
$ws.Range("M74").Value = -33.16669999999999

$ws.Range("H77").Value = 907.1667
$ws.Range("I77").Value = 907.1667
$ws.Range("K77").Value = 4535.8335
$ws.Range("M77").Value = -167.8334999999997

$ws.Range("H102").Value = 2388.5557
$ws.Range("I102").Value = 2642.5715
$ws.Range("J102").Value = 1499.5
$ws.Range("K102").Value = 2642.5715
$ws.Range("L102").Value = 1499.5
$ws.Range("M102").Value = -1020.5715
$ws.Range("N102").Value = -4743.5

$ws.Range("H116").Value = 3096.8572
$ws.Range("I116").Value = 1056.4445
$ws.Range("K116").Value = 1056.4445
$ws.Range("M116").Value = 1237.5555

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3096.8572
$ws.Range("I3").Value = 1056.4445
$ws.Range("K3").Value = 1056.4445
$ws.Range("M3").Value = -942.4445000000001

$ws.Range("H25").Value = 1333.3334
$ws.Range("I25").Value = 1333.3334
$ws.Range("K25").Value = 1333.3334
$ws.Range("M25").Value = -1098.3334

$ws.Range("H94").Value = 4833.3335
$ws.Range("J94").Value = 5000
$ws.Range("L94").Value = 5000
$ws.Range("N94").Value = -5902

$ws.Range("H134").Value = 1756.0714
$ws.Range("I134").Value = 1660.3846
$ws.Range("K134").Value = 4981.1538
$ws.Range("M134").Value = -2446.1538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2013
$ws.Range("J16").Value = 2013
$ws.Range("L16").Value = 2013
$ws.Range("N16").Value = -2587

$ws.Range("H58").Value = 1065.6666
$ws.Range("I58").Value = 1065.6666
$ws.Range("K58").Value = 1065.6666
$ws.Range("M58").Value = -862.6666

$ws.Range("H69").Value = 2166.6667
$ws.Range("I69").Value = 2166.6667
$ws.Range("K69").Value = 2166.6667
$ws.Range("M69").Value = -1417.6667

$ws.Range("H72").Value = 2166.6667
$ws.Range("I72").Value = 2166.6667
$ws.Range("K72").Value = 6500.000100000001
$ws.Range("M72").Value = -2756.000100000001

$ws.Range("H99").Value = 2988.889
$ws.Range("I99").Value = 3733.3333
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 3733.3333
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -2235.3333
$ws.Range("N99").Value = -4496

$ws.Range("H113").Value = 2013
$ws.Range("J113").Value = 2013
$ws.Range("L113").Value = 2013
$ws.Range("N113").Value = -6353

$ws.Range("H126").Value = 2988.889
$ws.Range("I126").Value = 3733.3333
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 11199.9999
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -8729.999899999999
$ws.Range("N126").Value = -9440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 2138
$ws.Range("I134").Value = 1700
$ws.Range("K134").Value = 5100
$ws.Range("M134").Value = -2565

$ws.Range("H136").Value = 1065.6666
$ws.Range("I136").Value = 1065.6666
$ws.Range("K136").Value = 3196.9998
$ws.Range("M136").Value = -646.9998000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H131").Value = 1003.1
$ws.Range("I131").Value = 786.625
$ws.Range("J131").Value = 1081.8182
$ws.Range("K131").Value = 2359.875
$ws.Range("L131").Value = 3245.4546
$ws.Range("M131").Value = 2680.125
$ws.Range("N131").Value = -13325.4546

$ws.Range("H24").Value = 42592.168
$ws.Range("J24").Value = 42592.168
$ws.Range("L24").Value = 42592.168
$ws.Range("N24").Value = -42938.168

$ws.Range("H70").Value = 8900.200000000001
$ws.Range("I70").Value = 8900.200000000001
$ws.Range("K70").Value = 8900.200000000001
$ws.Range("M70").Value = -8630.200000000001

$ws.Range("H73").Value = 8900.200000000001
$ws.Range("I73").Value = 8900.200000000001
$ws.Range("K73").Value = 8900.200000000001
$ws.Range("M73").Value = -7964.200000000001

$ws.Range("H80").Value = 2861.2
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 3326.5
$ws.Range("K80").Value = 1000
$ws.Range("L80").Value = 3326.5
$ws.Range("M80").Value = -2
$ws.Range("N80").Value = -5322.5

$ws.Range("H83").Value = 2861.2
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 3326.5
$ws.Range("K83").Value = 5000
$ws.Range("L83").Value = 16632.5
$ws.Range("M83").Value = -8
$ws.Range("N83").Value = -26616.5

$ws.Range("H97").Value = 1750.8182
$ws.Range("I97").Value = 1306.5555
$ws.Range("J97").Value = 3750
$ws.Range("K97").Value = 1306.5555
$ws.Range("L97").Value = 3750
$ws.Range("M97").Value = -810.5554999999999
$ws.Range("N97").Value = -4742

$ws.Range("H102").Value = 2049.6667
$ws.Range("I102").Value = 1964.3529
$ws.Range("K102").Value = 1964.3529
$ws.Range("M102").Value = -342.3529000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 39944
$ws.Range("J38").Value = 39944
$ws.Range("L38").Value = 39944
$ws.Range("N38").Value = -40764

$ws.Range("H40").Value = 2674.5
$ws.Range("I40").Value = 2674.5
$ws.Range("K40").Value = 2674.5
$ws.Range("M40").Value = -2538.5

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = $null

$ws.Range("H132").Value = 6256.6665
$ws.Range("I132").Value = 6256.6665
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 18769.9995
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -16239.9995
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 2517500
$ws.Range("J18").Value = 2517500
$ws.Range("L18").Value = 2517500
$ws.Range("N18").Value = -2517846

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null

$ws.Range("H107").Value = 988.1818
$ws.Range("I107").Value = 1009.125
$ws.Range("J107").Value = 932.3333
$ws.Range("K107").Value = 3027.375
$ws.Range("L107").Value = 2796.9999
$ws.Range("M107").Value = -1107.375
$ws.Range("N107").Value = -6636.9999

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = $null
